$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix column definitions: column A alone should carry the 30.71 width/style 1 ---
# (previously columns A:B shared that definition; now only column A does)
$ws.Columns.Item(1).ColumnWidth = 30.7109375

# --- New / updated shared text values used below ---
$objetivosPt = "Familiarizar o aluno com os conceitos básicos de equações diferenciais e suas aplicações."
$docente = "6270264 - Juan Fernando Zapata Zapata"
$programaResumidoPt = "Sequencias e séries, equações diferenciais ordinárias de 1ª e 2ª ordem com aplicações, solução de equações diferenciais por series de potencia, Séries de Fourier e Problemas de valores de contorno."
$programaPt = "Sequências e séries: Critérios de convergência, convergência condicional e absoluta, séries de potência, raio de convergência, derivação e integração termo a termo. Equações diferenciais ordinárias de 1ª e 2ª ordem: Equações exatas e não exatas, redução de ordem, Equação de Bernulli, método de variação de parâmetros e coeficientes a determinar, solução por séries de potencia de equações diferenciais, aplicações das equações diferenciais de 1ª e 2ª ordem.•Séries de Fourier: Teorema de convergência das séries de Fourier, Desigualdade de Bessel e Identidade de Parseval, equações em derivadas parciais e problemas de valores de contorno."
$metodoVal = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$criterioVal = "NF≥ 5,0."
$normaRecVal = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$bibliografiaVal = "1.H. L. Guidorizzi, UM CURSO DE CÁLCULO, volume IV. Livros Técnicos e Científicos, 1987.2.BRANNAN, James R. BOYCE, W.E. Equações diferenciais: uma Introdução a métodos modernos e suas aplicações. Rio de Janeiro: LTC ED., 2008.3.ZILL, D.G. ; CULLEN, M.R. Equações Diferenciais São Paulo: Pearson Makron Books2006., v.1 e 2.4.W. Kaplan, CÁLCULO AVANÇADO, volume II, Edgard Blücher, São Paulo, 1972.5.BOYCE,W.E. ; DIPRIMA,R.C. Equações diferenciais e problemas de valores de contorno. 8.ed. Rio de Janeiro: LTC Editora, 2008."
$requisito1 = "LOB1004 -  Cálculo II  (Requisito fraco)`n"
$requisito2 = "LOB1037 -  Àlgebra Linear  (Requisito fraco)`n"

# Row 10 ("Objetivos:") previously held the docente info by mistake; now holds the real PT objectives text
$ws.Range("B10").Value = $objetivosPt
$ws.Range("C10").Value = $objetivosPt

# --- Insert a new row at 13 (pushes everything from old row 13 onward down by one) ---
$ws.Rows.Item(13).Insert()

# Row 13 used to hold "Programa resumido:/Semestral" but that whole block moved down;
# the new row 13 now only carries the docente responsavel info (A13 stays empty).
$ws.Range("A13").Clear()
$ws.Range("B13").Value = $docente
$ws.Range("C13").Value = $docente
# The inserted row defaults B13 to the bold/label style; restore it to the normal wrap style
$ws.Range("B13").Font.Bold = $false
$ws.Range("B13").WrapText = $true
$ws.Range("B13").VerticalAlignment = -4160

# Row 14 (old row 13, "Programa resumido:/ Semestral") gets the new Portuguese summary text
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = $programaResumidoPt
$ws.Range("C14").Value = $programaResumidoPt
$ws.Rows.Item(14).RowHeight = 60

# Row 15 (old row 14 "Short syllabus:") keeps its English short-syllabus text, unchanged
$ws.Rows.Item(15).RowHeight = 60

# Row 16 (old row 15 "Programa:") previously mistakenly held the date; now holds full PT program text
$ws.Range("B16").Value = $programaPt
$ws.Range("C16").Value = $programaPt
$ws.Rows.Item(16).RowHeight = 120

# Row 17 (old row 16 "Syllabus:") keeps its English syllabus text, unchanged
$ws.Rows.Item(17).RowHeight = 120

# Row 18 (old row 17 "Avaliação:") stays label-only (no explicit custom height)

# Row 19 (old row 18 "Método:") previously mistakenly held docente info; now holds method/evaluation text
$ws.Range("B19").Value = $metodoVal
$ws.Range("C19").Value = $metodoVal
$ws.Rows.Item(19).RowHeight = 60

# Row 20 (old row 19 "Critério:") previously held the method text; now holds the criterio text
$ws.Range("B20").Value = $criterioVal
$ws.Range("C20").Value = $criterioVal
$ws.Rows.Item(20).RowHeight = 60

# Row 21 (old row 20 "Norma de recuperação:") previously held criterio text; now holds norma-recuperacao text
$ws.Range("B21").Value = $normaRecVal
$ws.Range("C21").Value = $normaRecVal
$ws.Rows.Item(21).RowHeight = 60  # was 120 before the edit, now 60

# Row 22 (old row 21 "Bibliografia:") previously held norma-recuperacao text; now holds full bibliografia text
$ws.Range("B22").Value = $bibliografiaVal
$ws.Range("C22").Value = $bibliografiaVal
$ws.Rows.Item(22).RowHeight = 120

# Row 23 (old row 22 "Requisitos:") stays label-only (no explicit custom height)

# Row 24 (old row 23, first requisito) unchanged content
$ws.Rows.Item(24).RowHeight = 30

# Row 25 (old row 24, second requisito) unchanged content
$ws.Rows.Item(25).RowHeight = 30
